$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.277.76'
$ws.Range('E2').Value = '  +4.04%  '

$ws.Range('D3').Value = '1.785.46'
$ws.Range('E3').Value = '  +0.06%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.14%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '339.17'
$ws.Range('E5').Value = '  +0.51%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9998'

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3819'
$ws.Range('E7').Value = '  -1.34%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3446'
$ws.Range('E8').Value = '  +0.49%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '47.06'
$ws.Range('E9').Value = '  -1.76%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.154'
$ws.Range('E10').Value = '  -2.91%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07400'
$ws.Range('E11').Value = '  -0.56%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '23.22'
$ws.Range('E12').Value = '  +7.20%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.001'
$ws.Range('E13').Value = '  -0.08%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.458'
$ws.Range('E14').Value = '  +0.49%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.342'
$ws.Range('E15').Value = '  +3.18%  '

$ws.Range('D16').Value = '1.780.99'
$ws.Range('E16').Value = '  -0.12%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001077'
$ws.Range('E17').Value = '  -1.21%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06672'
$ws.Range('E18').Value = '  +0.34%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '82.53'
$ws.Range('E19').Value = '  -0.90%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9998'
$ws.Range('E20').Value = '  -0.02%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.59'
$ws.Range('E21').Value = '  +0.08%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.466'
$ws.Range('E22').Value = '  -0.66%  '

$ws.Range('D23').Value = '28.256.93'
$ws.Range('E23').Value = '  +3.96%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.07'
$ws.Range('E24').Value = '  -2.15%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.375'
$ws.Range('E25').Value = '  +0.72%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.457'
$ws.Range('E26').Value = '  +0.79%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.74'
$ws.Range('E27').Value = '  -1.80%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.426'
$ws.Range('E28').Value = '  -2.66%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '153.82'
$ws.Range('E29').Value = '  -1.76%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '136.59'
$ws.Range('E30').Value = '  +1.98%  '

$ws.Range('D31').Value = '1.981.54'
$ws.Range('E31').Value = '  -0.22%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.143'
$ws.Range('E32').Value = '  +2.83%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.933'
$ws.Range('E33').Value = '  -1.07%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08876'

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '12.81'
$ws.Range('E35').Value = '  -0.90%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02446'
$ws.Range('E36').Value = '  +4.46%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6845'
$ws.Range('E37').Value = '  +0.59%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.329'
$ws.Range('E38').Value = '  -1.25%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06384'
$ws.Range('E39').Value = '  +0.77%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2178'
$ws.Range('E40').Value = '  -0.40%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.245'
$ws.Range('E41').Value = '  +0.45%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.499'
$ws.Range('E42').Value = '  -7.48%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.315'
$ws.Range('E43').Value = '  -1.49%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.13'
$ws.Range('E44').Value = '  -1.04%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9993'
$ws.Range('E45').Value = '  -0.03%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6308'
$ws.Range('E46').Value = '  -1.23%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.878'
$ws.Range('E47').Value = '  +0.52%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '133.51'
$ws.Range('E48').Value = '  +1.62%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.096'
$ws.Range('E49').Value = '  -2.21%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07452'
$ws.Range('E50').Value = '  +4.66%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.207'
$ws.Range('E51').Value = '  +8.12%  '
